$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the previously-missing AA1 cell and correct the sequence that
# follows it (AB1/AC1 shift from 25/26 to 26/27).
$ws.Range("AA1").Value = 25
$ws.Range("AB1").Value = 26
$ws.Range("AC1").Value = 27

# Board rows 12 & 13: mark a couple of cells with the new "a"/"b" tokens.
$ws.Range("C12").Value = "a"
$ws.Range("AB12").Value = "b"
$ws.Range("C13").Value = "a"
$ws.Range("AB13").Value = "b"

# Column AA was manually widened (losing its "best fit" autosize flag).
$ws.Columns.Item(27).ColumnWidth = 3.3072916666666665

# Move the active selection to AC11.
$null = $ws.Range("AC11").Select()
